$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 14.81967709108768
$ws.Range("C2").Value = 8.701991038463163
$ws.Range("D2").Value = 7.740546356698323
$ws.Range("E2").Value = 13.11470004589203
$ws.Range("F2").Value = 39.8029137826646
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 10.43965123685456
$ws.Range("K2").Value = 10.97304233228029
$ws.Range("L2").Value = 10.64414207485449
$ws.Range("O2").Value = 30.85463545647012
$ws.Range("B3").Value = 14.62931646709674
$ws.Range("C3").Value = 8.691439001280012
$ws.Range("D3").Value = 7.723530605890964
$ws.Range("E3").Value = 13.13361988409838
$ws.Range("F3").Value = 39.89908891827515
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 10.4610844188764
$ws.Range("K3").Value = 10.83240246286643
$ws.Range("L3").Value = 10.64266630505119
$ws.Range("O3").Value = 30.94776686950027
$ws.Range("B4").Value = 14.51376867355322
$ws.Range("C4").Value = 8.685101927765661
$ws.Range("D4").Value = 7.714087536435392
$ws.Range("E4").Value = 13.14684489049911
$ws.Range("F4").Value = 39.96572229615061
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 10.47505103561537
$ws.Range("K4").Value = 10.74678159678719
$ws.Range("L4").Value = 10.64300740739914
$ws.Range("O4").Value = 31.01030512843781
$ws.Range("B5").Value = 14.46707100463806
$ws.Range("C5").Value = 8.682555940822414
$ws.Range("D5").Value = 7.710494774059352
$ws.Range("E5").Value = 13.15263905169099
$ws.Range("F5").Value = 39.99478004119007
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 10.48094584146439
$ws.Range("K5").Value = 10.71211283334712
$ws.Range("L5").Value = 10.64346122514577
$ws.Range("O5").Value = 31.03713545824352
$ws.Range("B6").Value = 14.4593419383035
$ws.Range("C6").Value = 8.682135402780364
$ws.Range("D6").Value = 7.709913698467634
$ws.Range("E6").Value = 13.15362563393736
$ws.Range("F6").Value = 39.99971998852551
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 10.48193696224478
$ws.Range("K6").Value = 10.7063706581296
$ws.Range("L6").Value = 10.64355562838786
$ws.Range("O6").Value = 31.04167184409232
$ws.Range("B7").Value = 14.51313724720179
$ws.Range("C7").Value = 8.685067443256928
$ws.Range("D7").Value = 7.714038045733258
$ws.Range("E7").Value = 13.14692139268823
$ws.Range("F7").Value = 39.96610647267224
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 10.47512971125002
$ws.Range("K7").Value = 10.74631309077746
$ws.Range("L7").Value = 10.64301225155649
$ws.Range("O7").Value = 31.01066152538205
$ws.Range("B8").Value = 14.75379599158113
$ws.Range("C8").Value = 8.698323681062561
$ws.Range("D8").Value = 7.734472495406212
$ws.Range("E8").Value = 13.12089013287824
$ws.Range("F8").Value = 39.834500111844
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 10.44687428678337
$ws.Range("K8").Value = 10.92442064385762
$ws.Range("L8").Value = 10.6433751035986
$ws.Range("O8").Value = 30.88563503520244
$ws.Range("B9").Value = 15.23395264602104
$ws.Range("C9").Value = 8.725413719840896
$ws.Range("D9").Value = 7.782389371023138
$ws.Range("E9").Value = 13.08258103744324
$ws.Range("F9").Value = 39.63667211046077
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 10.39784338144481
$ws.Range("K9").Value = 11.27781867224495
$ws.Range("L9").Value = 10.65392819927299
$ws.Range("O9").Value = 30.68299857646559
$ws.Range("B10").Value = 15.58858143720893
$ws.Range("C10").Value = 8.745943010902064
$ws.Range("D10").Value = 7.822197323066357
$ws.Range("E10").Value = 13.06217048981788
$ws.Range("F10").Value = 39.5281792940071
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 10.36567830589275
$ws.Range("K10").Value = 11.53773602245048
$ws.Range("L10").Value = 10.6676042628551
$ws.Range("O10").Value = 30.56012022536082
$ws.Range("B11").Value = 15.74965573412404
$ws.Range("C11").Value = 8.755409480334276
$ws.Range("D11").Value = 7.841265088907927
$ws.Range("E11").Value = 13.05455763580052
$ws.Range("F11").Value = 39.48684563244934
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 10.35187691088282
$ws.Range("K11").Value = 11.65557283967953
$ws.Range("L11").Value = 10.67509371608225
$ws.Range("O11").Value = 30.50987806890997
$ws.Range("B12").Value = 15.81056100720922
$ws.Range("C12").Value = 8.759011739754831
$ws.Range("D12").Value = 7.848619669007608
$ws.Range("E12").Value = 13.05191459984156
$ws.Range("F12").Value = 39.47234814513764
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 10.34676964745894
$ws.Range("K12").Value = 11.70009901126915
$ws.Range("L12").Value = 10.67811041067315
$ws.Range("O12").Value = 30.49166699845185
$ws.Range("B13").Value = 15.7974488231034
$ws.Range("C13").Value = 8.758235163221215
$ws.Range("D13").Value = 7.847029830576124
$ws.Range("E13").Value = 13.05247317111006
$ws.Range("F13").Value = 39.4754190655425
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 10.34786430101244
$ws.Range("K13").Value = 11.6905143767688
$ws.Range("L13").Value = 10.67745270820885
$ws.Range("O13").Value = 30.49555283219066
$ws.Range("B14").Value = 15.75466855903177
$ws.Range("C14").Value = 8.755705494209321
$ws.Range("D14").Value = 7.841867491098629
$ws.Range("E14").Value = 13.05433538981727
$ws.Range("F14").Value = 39.48562976733916
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 10.35145435020533
$ws.Range("K14").Value = 11.65923817634425
$ws.Range("L14").Value = 10.67533829224245
$ws.Range("O14").Value = 30.50836350163529
$ws.Range("B15").Value = 15.72845109311715
$ws.Range("C15").Value = 8.754158254684649
$ws.Range("D15").Value = 7.83872274678234
$ws.Range("E15").Value = 13.05550725991825
$ws.Range("F15").Value = 39.49203452487144
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 10.35366884714982
$ws.Range("K15").Value = 11.64006696670712
$ws.Range("L15").Value = 10.67406661726952
$ws.Range("O15").Value = 30.51631651994649
$ws.Range("B16").Value = 15.57804517883591
$ws.Range("C16").Value = 8.745326851171759
$ws.Range("D16").Value = 7.820970169095419
$ws.Range("E16").Value = 13.06270159857961
$ws.Range("F16").Value = 39.53104205833036
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 10.36659693828437
$ws.Range("K16").Value = 11.53002373369565
$ws.Range("L16").Value = 10.66714017980721
$ws.Range("O16").Value = 30.56351758589155
$ws.Range("B17").Value = 15.48567333709915
$ws.Range("C17").Value = 8.73994115923308
$ws.Range("D17").Value = 7.810322442129934
$ws.Range("E17").Value = 13.06754288780133
$ws.Range("F17").Value = 39.55702717098831
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 10.37474035316442
$ws.Range("K17").Value = 11.46238543053439
$ws.Range("L17").Value = 10.66321461743764
$ws.Range("O17").Value = 30.59392340894325
$ws.Range("B18").Value = 15.43252331576842
$ws.Range("C18").Value = 8.736855511680428
$ws.Range("D18").Value = 7.804288650900074
$ws.Range("E18").Value = 13.0704848879044
$ws.Range("F18").Value = 39.57272796320036
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 10.37950244055493
$ws.Range("K18").Value = 11.42344617205877
$ws.Range("L18").Value = 10.66107621259331
$ws.Range("O18").Value = 30.61194431663897
$ws.Range("B19").Value = 15.41452584573234
$ws.Range("C19").Value = 8.735812861053589
$ws.Range("D19").Value = 7.802261371494835
$ws.Range("E19").Value = 13.07150805433426
$ws.Range("F19").Value = 39.57817357987943
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 10.38112824720555
$ws.Range("K19").Value = 11.41025709967085
$ws.Range("L19").Value = 10.66037275731486
$ws.Range("O19").Value = 30.61813726739933
$ws.Range("B20").Value = 15.49550896449618
$ws.Range("C20").Value = 8.740513233190699
$ws.Range("D20").Value = 7.811446571974935
$ws.Range("E20").Value = 13.06701123660771
$ws.Range("F20").Value = 39.55418287707854
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 10.37386538132522
$ws.Range("K20").Value = 11.4695895987158
$ws.Range("L20").Value = 10.66362014748202
$ws.Range("O20").Value = 30.59063155868898
$ws.Range("B21").Value = 15.76723703022528
$ws.Range("C21").Value = 8.75644805071498
$ws.Range("D21").Value = 7.843380189280738
$ws.Range("E21").Value = 13.05378190866978
$ws.Range("F21").Value = 39.4825992898456
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 10.35039663876375
$ws.Range("K21").Value = 11.66842766286208
$ws.Range("L21").Value = 10.67595446046982
$ws.Range("O21").Value = 30.5045785814305
$ws.Range("B22").Value = 15.94427789621936
$ws.Range("C22").Value = 8.766964060888302
$ws.Range("D22").Value = 7.865029996598991
$ws.Range("E22").Value = 13.04653310692512
$ws.Range("F22").Value = 39.44254600982487
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 10.3357520367963
$ws.Range("K22").Value = 11.79780280967126
$ws.Range("L22").Value = 10.68506744734967
$ws.Range("O22").Value = 30.45308629831329
$ws.Range("B23").Value = 15.84985568064157
$ws.Range("C23").Value = 8.76134244189519
$ws.Range("D23").Value = 7.853405090310131
$ws.Range("E23").Value = 13.05027429805415
$ws.Range("F23").Value = 39.46330696550741
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 10.34350481243198
$ws.Range("K23").Value = 11.72881806572087
$ws.Range("L23").Value = 10.68010803968924
$ws.Range("O23").Value = 30.4801338352039
$ws.Range("B24").Value = 15.49106241092093
$ws.Range("C24").Value = 8.740254565294267
$ws.Range("D24").Value = 7.810938079328677
$ws.Range("E24").Value = 13.067251101804
$ws.Range("F24").Value = 39.55546641009489
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 10.37426070584053
$ws.Range("K24").Value = 11.46633275672642
$ws.Range("L24").Value = 10.663436438253
$ws.Range("O24").Value = 30.59211812158613
$ws.Range("B25").Value = 15.10350454338539
$ws.Range("C25").Value = 8.717972166717331
$ws.Range("D25").Value = 7.768604543889522
$ws.Range("E25").Value = 13.09158393568018
$ws.Range("F25").Value = 39.68372597824214
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 10.41042788225338
$ws.Range("K25").Value = 11.18200601018616
$ws.Range("L25").Value = 10.65002727404116
$ws.Range("O25").Value = 30.73325596472722
